# working fast ref pulse output for ADC alignment
#
# Adds chip-ID registers (rows 11-13), RF switch/pulse control registers
# (rows 48-49), a stray space label at C86, and moves the active
# selection/scroll position back up to C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 86: stray blank label ---
$ws.Range("C86").Value = " "

# --- rows 48-49: new RF switch / RF pulse registers ---
$ws.Range("B48").Value = "RF switch select"
$ws.Range("B49").Value = "RF pulse enable"

$ws.Range("E48").Value = "pick off debug header?"
$ws.Range("E49").Value = "use trig_out aux SMA"

$ws.Range("C48").Value = "toggle input switch between signal and ref pulse (LSB=1 switch to cal pulse input)"
$ws.Range("C49").Value = "toggle FPGA-generated fast pulse (LSB=1 enable, LSB=0 disable)"

$ws.Range("D48").Value = "0x000000"
$ws.Range("D49").Value = "0x000000"

# --- rows 11-13: chip_id registers (replaces old status_1/status_2 block) ---
$ws.Range("B11").Value = "chip_id(low)"
$ws.Range("B12").Value = "chip_id(mid)"
$ws.Range("B13").Value = "chip_id(high)"

$ws.Range("E12").Value = "mid 24 bits"
$ws.Range("E11").Value = "chip is 64 bits: lower 24 bits"
$ws.Range("E13").Value = "high 16 bits"

$ws.Range("C13").Value = "read_only"

# --- restore view: scroll to top, select C10 ---
$ws.Range("C10").Select()
